# franklin-201602558 finishing data analysis
#
# Updates the "Val1" (column E) and "Val2" (column G) statistics for every
# material row on the active sheet. The source values are numeric-looking
# text (they were stored as text in the workbook), so we re-enter them with
# a leading apostrophe - exactly like a user typing a quoted numeric value
# into Excel - to keep them stored as text rather than being auto-converted
# to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (Val1 in column E, Val2 in column G)
$updates = @{
    2 = @("14.2227408704967", "1053.61913665593")
    3 = @("893.647058823529", "866.158282042129")
    4 = @("10.8034562186542", "883.73441756635")
    5 = @("746.5",            "83.1954643842899")
    6 = @("1130.37837837838", "1118.71435262794")
    7 = @("952.041666666667", "1178.55676143772")
    8 = @("18.5419634233197", "825.514186196713")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("E$row").Formula = "'" + $vals[0]
    $ws.Range("G$row").Formula = "'" + $vals[1]
}
